$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
$statusOld = "Ready for handoff"
$statusNew = "Handed back: in sync with en-US"

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("E3").Value = $statusNew
$wsOverview.Range("F3").Value = $statusNew

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsZh.Range("C2").Value = $statusNew
$wsZh.Range("C3").Value = $statusNew
$wsDe.Range("C2").Value = $statusNew
$wsDe.Range("C3").Value = $statusNew

# --- zh-cn: fill in "Latest Target File" (I) / "Latest Handback File" (J) / "Latest Handback DateTime" (K) ---
$wsZh.Range("I2").Value = $wsZh.Range("A2").Value
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/435bec39117bfc05cf1894f0da9963a51ae92c99/e2e/006ef815-8d7e-446b-8fa4-878927694f7c.md", "", "", "006ef815-8d7e-446b-8fa4-878927694f7c.md")
$wsZh.Range("J2").Value = $wsZh.Range("G2").Value
$wsZh.Range("K2").Value = "2016-08-29 10:28:02"

$wsZh.Range("I3").Value = $wsZh.Range("A3").Value
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/435bec39117bfc05cf1894f0da9963a51ae92c99/e2e/fe20e06e-5122-4407-91fe-e8d72a716447.md", "", "", "fe20e06e-5122-4407-91fe-e8d72a716447.md")
$wsZh.Range("J3").Value = $wsZh.Range("G3").Value
$wsZh.Range("K3").Value = "2016-08-29 10:28:02"

# --- de-de: fill in "Latest Target File" (I) / "Latest Handback File" (J) / "Latest Handback DateTime" (K) ---
$wsDe.Range("I2").Value = $wsDe.Range("A2").Value
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/435bec39117bfc05cf1894f0da9963a51ae92c99/e2e/006ef815-8d7e-446b-8fa4-878927694f7c.md", "", "", "006ef815-8d7e-446b-8fa4-878927694f7c.md")
$wsDe.Range("J2").Value = $wsDe.Range("G2").Value
$wsDe.Range("K2").Value = "2016-08-29 10:28:14"

$wsDe.Range("I3").Value = $wsDe.Range("A3").Value
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/435bec39117bfc05cf1894f0da9963a51ae92c99/e2e/fe20e06e-5122-4407-91fe-e8d72a716447.md", "", "", "fe20e06e-5122-4407-91fe-e8d72a716447.md")
$wsDe.Range("J3").Value = $wsDe.Range("G3").Value
$wsDe.Range("K3").Value = "2016-08-29 10:28:14"

# --- Column width adjustments to fit the newly-populated hyperlink columns ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

Write-Host "Report regenerated for handback."
